# Reporte_juego.xlsx — "enviar cambios de optimizacion"
#
# Adds a new "Juego" column (X) that flags whether a played row came from a
# simulation or not, and appends the most recent play as a new row (36).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New header cell X1 = "Juego", matching the style of the other
#        header cells (bold, centered, bordered) by copying W1's format. ---
$ws.Range("W1").Copy() | Out-Null
$ws.Range("X1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("X1").Value = "Juego"

# --- 2. Existing rows (2-35) get a blank placeholder in the new column. ---
for ($r = 2; $r -le 35; $r++) {
    $ws.Cells.Item($r, 24).FormulaR1C1 = '=""'
}

# --- 3. Append the new play as row 36. ---
$row = 36
$ws.Cells.Item($row, 1).Value  = "2024-03-30 19:24:12"
$ws.Cells.Item($row, 2).Value  = 1
$ws.Cells.Item($row, 3).Value  = 1
$ws.Cells.Item($row, 4).Value  = 0
$ws.Cells.Item($row, 5).Value  = 1
$ws.Cells.Item($row, 6).Value  = 0
$ws.Cells.Item($row, 7).Value  = 0
$ws.Cells.Item($row, 8).Value  = 0
$ws.Cells.Item($row, 9).Value  = 0.001
$ws.Cells.Item($row, 10).Value = 0.05
$ws.Cells.Item($row, 11).Value = 0.003
$ws.Cells.Item($row, 12).Value = 100
$ws.Cells.Item($row, 13).Value = 500
$ws.Cells.Item($row, 14).Value = 10
$ws.Cells.Item($row, 15).Value = 9
$ws.Cells.Item($row, 16).Value = 2
$ws.Cells.Item($row, 17).Value = 1000
$ws.Cells.Item($row, 18).Value = 3
$ws.Cells.Item($row, 19).Value = 1
$ws.Cells.Item($row, 20).Value = 100
$ws.Cells.Item($row, 21).Value = 1
$ws.Cells.Item($row, 22).Value = "Data/bombay1.xlsx"
$ws.Cells.Item($row, 23).Value = 29000
$ws.Cells.Item($row, 24).Value = "No es Simulación"
